# Auto-generated script to update cryptos worksheet values
# per commit: "Updated cryptos list on Wed Sep  6 03:20:26 UTC 2023 with GitHub Actions"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The "Price" column (D) contains values that look numeric (e.g. "215.80",
# "0.5079") but must stay as literal text, exactly as authored in the
# source feed (trailing zeros, thousand-dot-grouping, subscript digits,
# etc). Force the column's number format to Text *before* writing any
# values so Excel does not silently convert them to floating point
# numbers (which would also corrupt formatting like trailing zeros).
$ws.Range('D2:D51').NumberFormat = '@'


$ws.Range('D2').Value = '25.900.31'
$ws.Range('E2').Value = '  +0.94%  '
$ws.Range('D3').Value = '1.641.89'
$ws.Range('E3').Value = '  +1.61%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = '215.80'
$ws.Range('E5').Value = '  +0.77%  '
$ws.Range('D6').Value = '0.5079'
$ws.Range('E6').Value = '  +0.45%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('E8').Value = '  +1.91%  '
$ws.Range('D9').Value = '0.06476'
$ws.Range('E9').Value = '  +2.10%  '
$ws.Range('D10').Value = '20.31'
$ws.Range('E10').Value = '  +5.73%  '
$ws.Range('D11').Value = '0.07811'
$ws.Range('E11').Value = '  +0.65%  '
$ws.Range('B12').Value = 'Polkadot'
$ws.Range('C12').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D12').Value = '4.273'
$ws.Range('E12').Value = '  +1.14%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.648.49'
$ws.Range('E13').Value = '  +1.82%  '
$ws.Range('D14').Value = '1.869.83'
$ws.Range('E14').Value = '  +1.52%  '
$ws.Range('D15').Value = '0.5667'
$ws.Range('E15').Value = '  +2.41%  '
$ws.Range('D16').Value = '0.0₅7726'
$ws.Range('E16').Value = '  +3.12%  '
$ws.Range('D17').Value = '63.63'
$ws.Range('E17').Value = '  +0.46%  '
$ws.Range('D18').Value = '25.929.70'
$ws.Range('E18').Value = '  +0.96%  '
$ws.Range('E19').Value = '  +0.10%  '
$ws.Range('D20').Value = '194.79'
$ws.Range('E20').Value = '  +0.81%  '
$ws.Range('D21').Value = '4.410'
$ws.Range('E21').Value = '  +1.64%  '
$ws.Range('E22').Value = '  +3.10%  '
$ws.Range('D23').Value = '6.283'
$ws.Range('E23').Value = '  +5.75%  '
$ws.Range('D24').Value = '1.003'
$ws.Range('E24').Value = '  +0.07%  '
$ws.Range('E25').Value = '  -4.07%  '
$ws.Range('D26').Value = '138.99'
$ws.Range('E26').Value = '  -0.78%  '
$ws.Range('D27').Value = '0.1230'
$ws.Range('E27').Value = '  -1.80%  '
$ws.Range('D28').Value = '6.873'
$ws.Range('E28').Value = '  +2.69%  '
$ws.Range('D29').Value = '15.64'
$ws.Range('E29').Value = '  +1.60%  '
$ws.Range('D30').Value = '1.246'
$ws.Range('E30').Value = '  +1.27%  '
$ws.Range('D31').Value = '0.05008'
$ws.Range('E31').Value = '  +3.58%  '
$ws.Range('D32').Value = '3.329'
$ws.Range('E32').Value = '  +1.54%  '
$ws.Range('D33').Value = '3.278'
$ws.Range('E33').Value = '  +3.70%  '
$ws.Range('E34').Value = '  +3.11%  '
$ws.Range('D35').Value = '2.384'
$ws.Range('E35').Value = '  +0.94%  '
$ws.Range('D36').Value = '0.9106'
$ws.Range('E36').Value = '  +2.59%  '
$ws.Range('D37').Value = '2.587'
$ws.Range('E37').Value = '  +2.43%  '
$ws.Range('D38').Value = '0.5547'
$ws.Range('E38').Value = '  +1.65%  '
$ws.Range('D39').Value = '1.130.65'
$ws.Range('E39').Value = '  +0.88%  '
$ws.Range('D40').Value = '0.01578'
$ws.Range('E40').Value = '  +1.56%  '
$ws.Range('D41').Value = '1.003'
$ws.Range('E41').Value = '  -2.48%  '
$ws.Range('D42').Value = '5.514'
$ws.Range('E42').Value = '  -0.66%  '
$ws.Range('D43').Value = '99.88'
$ws.Range('E43').Value = '  +3.17%  '
$ws.Range('D44').Value = '0.8018'
$ws.Range('E44').Value = '  +1.53%  '
$ws.Range('D45').Value = '0.0₈112'
$ws.Range('E45').Value = '  +0.65%  '
$ws.Range('D46').Value = '55.82'
$ws.Range('E46').Value = '  +2.59%  '
$ws.Range('D47').Value = '0.4235'
$ws.Range('E47').Value = '  -3.90%  '
$ws.Range('D48').Value = '7.708'
$ws.Range('E48').Value = '  +2.61%  '
$ws.Range('D49').Value = '0.05048'
$ws.Range('E49').Value = '  -0.85%  '
$ws.Range('D50').Value = '1.006'
$ws.Range('E50').Value = '  +1.22%  '
$ws.Range('D51').Value = '1.001'
$ws.Range('E51').Value = '  +0.03%  '
